# BBDC_inputs_bear.xlsx — re-saved/edited workbook (pt-BR Excel build)
#
# Content changes:
#  - rename Sheet1 -> Planilha1
#  - clear the "ATIVO" header text in A1
#  - recompute a handful of TARGET_BEAR-sheet figures (C3, C11, C12, C13)
#  - change the workbook default font size from 12 to 11
#  - update the selection to the full used range (A1:G13)
#  - page margins switch to the metric (cm) defaults
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Planilha1"

$ws.Range("A1").Value = ""

$ws.Range("C3").Value = 7.2887113573703308
$ws.Range("C11").Value = 4.9956599612547654
$ws.Range("C12").Value = 17.81122021073049
$ws.Range("C13").Value = 7.2887113573703308

$wb.Styles.Item("Normal").Font.Size = 11

$ws.PageSetup.LeftMargin = 36.850393728
$ws.PageSetup.RightMargin = 36.850393728
$ws.PageSetup.TopMargin = 56.692913399999995
$ws.PageSetup.BottomMargin = 56.692913399999995
$ws.PageSetup.HeaderMargin = 22.67716464
$ws.PageSetup.FooterMargin = 22.67716464

$ws.Range("A1:G13").Select() | Out-Null
